# Slide 4, TextBox "TextBox 388": split the run
#   " один до одного. Таким чином, усі колеса (в тому числі шестерні та
#   сотуари) мають однаковий розмір та незалежність від грузоподібних
#   машин. "
# into three runs:
#   " одна до "   (keeps the original "dirty" rPr: b="0" i="0" + empty effectLst)
#   "іншої"       (new run with "clean" rPr - no b/i/effectLst)
#   ". Таким чином, ... машин. "  (keeps the original "dirty" rPr)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)          # "TextBox 388"
$tr = $sh.TextFrame.TextRange

# The run immediately before the target text ("дотикаються") already has
# the "clean" rPr pattern we want for the new middle run. By replacing a
# span that STARTS inside that clean run, the engine reformats the whole
# replaced span using that clean formatting - this lets us manufacture a
# clean-formatted run instead of inheriting the "dirty" one.
$anchorAndTarget = $tr.Characters(499, 151)
$anchorAndTarget.Text = "дотикаються одна до іншої. Таким чином, усі колеса (в тому числі шестерні та сотуари) мають однаковий розмір та незалежність від грузоподібних машин. "

# Re-select (now all "clean") the two outer pieces and restore them back to
# the original "dirty" formatting (b="0" i="0" + empty effectLst), which
# also forces them to split off into their own runs, leaving "іншої" as an
# isolated, still-clean, middle run.
$before = $tr.Characters(510, 9)
$before.Text = " одна до "
$before.Font.Bold = 0
$before.Font.Italic = 0
$before.Font.Shadow = 0

$after = $tr.Characters(524, 125)
$after.Text = ". Таким чином, усі колеса (в тому числі шестерні та сотуари) мають однаковий розмір та незалежність від грузоподібних машин. "
$after.Font.Bold = 0
$after.Font.Italic = 0
$after.Font.Shadow = 0
